$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2025-09-25 Thursday" "2025-09-26 Friday"

# Multiplication table cells, applied in document order so that a
# replacement's new value never collides with an not-yet-processed
# old value.
Replace-Text "48×40=" "76×91="
Replace-Text "28×87=" "37×96="
Replace-Text "65×78=" "65×89="
Replace-Text "46×32=" "92×82="
Replace-Text "73×91=" "84×44="
Replace-Text "44×91=" "44×99="
Replace-Text "35×91=" "86×89="
Replace-Text "54×47=" "46×19="
Replace-Text "68×20=" "57×77="
Replace-Text "87×80=" "54×22="
Replace-Text "50×24=" "16×57="
Replace-Text "44×11=" "23×75="
Replace-Text "46×46=" "74×18="
Replace-Text "35×40=" "32×99="
Replace-Text "90×27=" "45×35="
Replace-Text "87×34=" "36×65="
Replace-Text "86×66=" "58×75="
Replace-Text "83×24=" "59×20="
Replace-Text "13×72=" "28×87="
Replace-Text "61×21=" "99×40="
Replace-Text "89×17=" "48×11="
Replace-Text "13×55=" "28×93="
Replace-Text "85×32=" "85×61="
Replace-Text "44×41=" "80×85="
Replace-Text "71×21=" "16×60="
